$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$helper = $ws.Range("ZZ1")

$ws.Range("D2").Value = '66.870.10'
$ws.Range("E2").Value = '  -3.95%  '
$ws.Range("D3").Value = '3.462.82'
$ws.Range("E3").Value = '  -4.25%  '
$helper.Value = "'1.00"
$helper.Copy()
$ws.Range("D4").PasteSpecial(-4104)
$ws.Range("E4").Value = '  +0.05%  '
$helper.Value = "'603.87"
$helper.Copy()
$ws.Range("D5").PasteSpecial(-4104)
$ws.Range("E5").Value = '  -4.30%  '
$helper.Value = "'147.47"
$helper.Copy()
$ws.Range("D6").PasteSpecial(-4104)
$ws.Range("E6").Value = '  -7.21%  '
$ws.Range("D7").Value = '3.461.11'
$helper.Value = "'0.999"
$helper.Copy()
$ws.Range("D8").PasteSpecial(-4104)
$ws.Range("E8").Value = '  -0.08%  '
$helper.Value = "'0.483"
$helper.Copy()
$ws.Range("D9").PasteSpecial(-4104)
$ws.Range("E9").Value = '  -2.39%  '
$ws.Range("E10").Value = '  -5.12%  '
$helper.Value = "'7.49"
$helper.Copy()
$ws.Range("D11").PasteSpecial(-4104)
$ws.Range("E11").Value = '  +0.12%  '
$helper.Value = "'0.423"
$helper.Copy()
$ws.Range("D12").PasteSpecial(-4104)
$ws.Range("E12").Value = '  -4.24%  '
$helper.Value = "'0.0000213"
$helper.Copy()
$ws.Range("D13").PasteSpecial(-4104)
$ws.Range("E13").Value = '  -6.54%  '
$helper.Value = "'31.68"
$helper.Copy()
$ws.Range("D14").PasteSpecial(-4104)
$ws.Range("E14").Value = '  -5.79%  '
$ws.Range("D15").Value = '4.041.91'
$ws.Range("E15").Value = '  -4.29%  '
$ws.Range("D16").Value = '3.472.69'
$ws.Range("E16").Value = '  -3.85%  '
$ws.Range("D17").Value = '66.858.98'
$ws.Range("E17").Value = '  -3.96%  '
$ws.Range("E18").Value = '  -0.89%  '
$helper.Value = "'6.42"
$helper.Copy()
$ws.Range("D19").PasteSpecial(-4104)
$ws.Range("E19").Value = '  -4.47%  '
$helper.Value = "'15.29"
$helper.Copy()
$ws.Range("D20").PasteSpecial(-4104)
$ws.Range("E20").Value = '  -5.59%  '
$helper.Value = "'10.00"
$helper.Copy()
$ws.Range("D21").PasteSpecial(-4104)
$ws.Range("E21").Value = '  -2.51%  '
$helper.Value = "'439.92"
$helper.Copy()
$ws.Range("D22").PasteSpecial(-4104)
$ws.Range("E22").Value = '  -4.96%  '
$helper.Value = "'0.607"
$helper.Copy()
$ws.Range("D23").PasteSpecial(-4104)
$ws.Range("E23").Value = '  -6.04%  '
$helper.Value = "'78.28"
$helper.Copy()
$ws.Range("D24").PasteSpecial(-4104)
$ws.Range("E24").Value = '  -0.67%  '
$ws.Range("E25").Value = '  +0.02%  '
$ws.Range("D26").Value = '3.595.75'
$ws.Range("E26").Value = '  -4.23%  '
$helper.Value = "'0.0000121"
$helper.Copy()
$ws.Range("D27").PasteSpecial(-4104)
$ws.Range("E27").Value = '  -10.78%  '
$ws.Range("E28").Value = '  -8.31%  '
$ws.Range("E29").Value = '  -9.54%  '
$ws.Range("E30").Value = '  -6.33%  '
$ws.Range("E31").Value = '  -7.01%  '
$ws.Range("B32").Value = 'Kaspa'
$ws.Range("C32").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$helper.Value = "'0.167"
$helper.Copy()
$ws.Range("D32").PasteSpecial(-4104)
$ws.Range("E32").Value = '  -4.46%  '
$ws.Range("B33").Value = 'Binance-PegBSC-USD'
$ws.Range("C33").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$helper.Value = "'1.00"
$helper.Copy()
$ws.Range("D33").PasteSpecial(-4104)
$ws.Range("E33").Value = '  +0.05%  '
$helper.Value = "'25.39"
$helper.Copy()
$ws.Range("D34").PasteSpecial(-4104)
$ws.Range("E34").Value = '  -4.49%  '
$helper.Value = "'6.10"
$helper.Copy()
$ws.Range("D35").PasteSpecial(-4104)
$ws.Range("E35").Value = '  -7.42%  '
$ws.Range("D36").Value = '3.451.69'
$ws.Range("E36").Value = '  -4.39%  '
$ws.Range("E37").Value = '  -7.83%  '
$ws.Range("B38").Value = 'Aptos'
$ws.Range("C38").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$helper.Value = "'7.92"
$helper.Copy()
$ws.Range("D38").PasteSpecial(-4104)
$ws.Range("E38").Value = '  -7.06%  '
$ws.Range("B39").Value = 'USDe'
$ws.Range("C39").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$helper.Value = "'1.00"
$helper.Copy()
$ws.Range("D39").PasteSpecial(-4104)
$ws.Range("E39").Value = '  +0.00%  '
$ws.Range("E40").Value = '  -0.06%  '
$helper.Value = "'173.56"
$helper.Copy()
$ws.Range("D41").PasteSpecial(-4104)
$ws.Range("E41").Value = '  -2.95%  '
$helper.Value = "'2.18"
$helper.Copy()
$ws.Range("D42").PasteSpecial(-4104)
$ws.Range("E42").Value = '  -9.50%  '
$helper.Value = "'0.0887"
$helper.Copy()
$ws.Range("D43").PasteSpecial(-4104)
$ws.Range("E43").Value = '  -4.67%  '
$helper.Value = "'5.37"
$helper.Copy()
$ws.Range("D44").PasteSpecial(-4104)
$ws.Range("E44").Value = '  -5.38%  '
$helper.Value = "'0.883"
$helper.Copy()
$ws.Range("D45").PasteSpecial(-4104)
$ws.Range("E45").Value = '  -3.42%  '
$helper.Value = "'28.88"
$helper.Copy()
$ws.Range("D46").PasteSpecial(-4104)
$ws.Range("E46").Value = '  -10.69%  '
$helper.Value = "'46.12"
$helper.Copy()
$ws.Range("D47").PasteSpecial(-4104)
$ws.Range("E47").Value = '  +0.32%  '
$helper.Value = "'1.24"
$helper.Copy()
$ws.Range("D48").PasteSpecial(-4104)
$ws.Range("E48").Value = '  -10.93%  '
$helper.Value = "'7.48"
$helper.Copy()
$ws.Range("D49").PasteSpecial(-4104)
$ws.Range("E49").Value = '  -4.60%  '
$helper.Value = "'2.46"
$helper.Copy()
$ws.Range("D50").PasteSpecial(-4104)
$ws.Range("E50").Value = '  -10.84%  '
$helper.Value = "'0.988"
$helper.Copy()
$ws.Range("D51").PasteSpecial(-4104)
$ws.Range("E51").Value = '  -4.55%  '

$helper.Clear()
$excel.CutCopyMode = $false
